$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5292.154
$ws.Range("I40").Value = 2999.8
$ws.Range("K40").Value = 2999.8
$ws.Range("M40").Value = -2824.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3377.5
$ws.Range("I62").Value = 2755
$ws.Range("K62").Value = 2755
$ws.Range("M62").Value = -2131

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3377.5
$ws.Range("I65").Value = 2755
$ws.Range("K65").Value = 13775
$ws.Range("M65").Value = -10655

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 3000
$ws.Range("J97").Value = 3000
$ws.Range("L97").Value = 9000
$ws.Range("N97").Value = -9992

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 935.62964
$ws.Range("I98").Value = 977.04346
$ws.Range("K98").Value = 977.04346
$ws.Range("M98").Value = 520.95654

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1911.4375
$ws.Range("I106").Value = 1892.2
$ws.Range("K106").Value = 1892.2
$ws.Range("M106").Value = -1261.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 935.62964
$ws.Range("I122").Value = 977.04346
$ws.Range("K122").Value = 2931.13038
$ws.Range("M122").Value = -481.1303800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1644.8833
$ws.Range("I32").Value = 1677.5344
$ws.Range("K32").Value = 1677.5344
$ws.Range("M32").Value = -1390.5344

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3551.0667
$ws.Range("I45").Value = 1956
$ws.Range("K45").Value = 1956
$ws.Range("M45").Value = -1579

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3047.4666
$ws.Range("I74").Value = 2268.3438
$ws.Range("K74").Value = 2268.3438
$ws.Range("M74").Value = -1394.3438

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3047.4666
$ws.Range("I77").Value = 2268.3438
$ws.Range("K77").Value = 11341.719
$ws.Range("M77").Value = -6973.719000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4987.769
$ws.Range("I132").Value = 4199.0312
$ws.Range("J132").Value = 8593.429
$ws.Range("K132").Value = 12597.0936
$ws.Range("L132").Value = 25780.287
$ws.Range("M132").Value = -10067.0936
$ws.Range("N132").Value = -30840.287

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4392.778
$ws.Range("I31").Value = 3329.7646
$ws.Range("K31").Value = 3329.7646
$ws.Range("M31").Value = -3034.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4392.778
$ws.Range("I34").Value = 3329.7646
$ws.Range("K34").Value = 3329.7646
$ws.Range("M34").Value = -3127.7646

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 22997.5
$ws.Range("J51").Value = 27996.666
$ws.Range("L51").Value = 27996.666
$ws.Range("N51").Value = -29468.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7382.3335
$ws.Range("I58").Value = 3756.5
$ws.Range("K58").Value = 3756.5
$ws.Range("M58").Value = -3553.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 22997.5
$ws.Range("J61").Value = 27996.666
$ws.Range("L61").Value = 27996.666
$ws.Range("N61").Value = -28692.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H95").Value = 37399.8
$ws.Range("J95").Value = 37399.8
$ws.Range("L95").Value = 37399.8
$ws.Range("N95").Value = -42891.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 19000
$ws.Range("J96").Value = 19000
$ws.Range("L96").Value = 19000
$ws.Range("N96").Value = -24492

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1923
$ws.Range("I105").Value = 1812.8
$ws.Range("K105").Value = 1812.8
$ws.Range("M105").Value = -65.79999999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3521.5356
$ws.Range("I122").Value = 3296.5
$ws.Range("K122").Value = 9889.5
$ws.Range("M122").Value = -7439.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7382.3335
$ws.Range("I136").Value = 3756.5
$ws.Range("K136").Value = 11269.5
$ws.Range("M136").Value = -8719.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 32254.445
$ws.Range("I141").Value = 30881.666
$ws.Range("K141").Value = 30881.666
$ws.Range("M141").Value = -25701.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 782.4286
$ws.Range("J38").Value = 1031.4
$ws.Range("L38").Value = 3094.2
$ws.Range("N38").Value = -3788.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 11304.833
$ws.Range("I126").Value = 7610
$ws.Range("K126").Value = 22830
$ws.Range("M126").Value = -17890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 23813176
$ws.Range("I131").Value = 250000750
$ws.Range("J131").Value = 3957.5789
$ws.Range("K131").Value = 750002250
$ws.Range("L131").Value = 11872.7367
$ws.Range("M131").Value = -749997210
$ws.Range("N131").Value = -21952.7367

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3836.8
$ws.Range("I102").Value = 3052.5715
$ws.Range("K102").Value = 3052.5715
$ws.Range("M102").Value = -1430.5715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4302
$ws.Range("I122").Value = 3856.75
$ws.Range("J122").Value = 6083
$ws.Range("K122").Value = 11570.25
$ws.Range("L122").Value = 18249
$ws.Range("M122").Value = -9120.25
$ws.Range("N122").Value = -23149

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5014.6665
$ws.Range("I126").Value = 4786.154
$ws.Range("J126").Value = 6500
$ws.Range("K126").Value = 14358.462
$ws.Range("L126").Value = 19500
$ws.Range("M126").Value = -11888.462
$ws.Range("N126").Value = -24440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H134").Value = 159284.86
$ws.Range("J134").Value = 159284.86
$ws.Range("L134").Value = 477854.58
$ws.Range("N134").Value = -482924.58

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 950
$ws.Range("J16").Value = 2300.6667
$ws.Range("L16").Value = 2300.6667
$ws.Range("N16").Value = -2640.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2759.125
$ws.Range("I40").Value = 2759.125
$ws.Range("K40").Value = 2759.125
$ws.Range("M40").Value = -2623.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2734.842
$ws.Range("I93").Value = 3240.889
$ws.Range("K93").Value = 3240.889
$ws.Range("M93").Value = -1992.889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6817.364
$ws.Range("I100").Value = 5874.5
$ws.Range("J100").Value = 7356.143
$ws.Range("K100").Value = 5874.5
$ws.Range("L100").Value = 7356.143
$ws.Range("M100").Value = -5333.5
$ws.Range("N100").Value = -8438.143

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4396.8
$ws.Range("I122").Value = 3985
$ws.Range("J122").Value = 4499.75
$ws.Range("K122").Value = 11955
$ws.Range("L122").Value = 13499.25
$ws.Range("M122").Value = -9505
$ws.Range("N122").Value = -18399.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 10333.333
$ws.Range("J43").Value = 11000
$ws.Range("L43").Value = 11000
$ws.Range("N43").Value = -11298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5063.1763
$ws.Range("I62").Value = 4876.8184
$ws.Range("J62").Value = 5404.8335
$ws.Range("K62").Value = 4876.8184
$ws.Range("L62").Value = 5404.8335
$ws.Range("M62").Value = -4252.8184
$ws.Range("N62").Value = -6652.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 5063.1763
$ws.Range("I65").Value = 4876.8184
$ws.Range("J65").Value = 5404.8335
$ws.Range("K65").Value = 24384.092
$ws.Range("L65").Value = 27024.1675
$ws.Range("M65").Value = -21264.092
$ws.Range("N65").Value = -33264.1675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1900
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2691.84
$ws.Range("I122").Value = 1853.0526
$ws.Range("K122").Value = 5559.1578
$ws.Range("M122").Value = -3109.1578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6504.5454
$ws.Range("I126").Value = 6255.75
$ws.Range("K126").Value = 18767.25
$ws.Range("M126").Value = -16297.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
